$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text format for Price column cells we touch that look numeric,
# so Excel keeps them stored as text (matching the source data).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.318.13"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.539.33"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.33"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.43"
$ws.Range("E6").Value = "  +2.63%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.536.66"
$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("E10").Value = "  +2.02%  "

$ws.Range("E11").Value = "  +2.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.95"
$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.998.38"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("E16").Value = "  +1.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.181.60"
$ws.Range("E17").Value = "  +2.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.494.12"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.53"
$ws.Range("E20").Value = "  +2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.60"
$ws.Range("E21").Value = "  +4.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.22"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.71"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.17"
$ws.Range("E26").Value = "  +2.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.24"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("E30").Value = "  +1.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.32"
$ws.Range("E31").Value = "  +2.49%  "

$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("E37").Value = "  +0.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.69"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.88"
$ws.Range("E39").Value = "  +1.43%  "

$ws.Range("E40").Value = "  +1.65%  "

$ws.Range("E41").Value = "  +0.99%  "

$ws.Range("E42").Value = "  +2.53%  "

$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.564"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.68"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("E49").Value = "  +2.03%  "

$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("E51").Value = "  +1.22%  "
